# "finish the ruintown event to test the next function"
#
# 1) The old row 24 (Id 42010015 "传送门/portal", a stray duplicate of the
#    existing portal-event row) is removed and re-created as a brand new
#    first entry (Id 42010001) inserted at row 11 - this fills a gap that
#    existed in the Id numbering (...42010002..42010016 were used, but
#    42010001 was never used).
# 2) Everything that used to sit at rows 11-23 shifts down one row to make
#    room (12-24).
# 3) The row that used to be the last one (old row 25, Id 42010016) stays
#    in place as row 25, but is renumbered to 42010015 now that the id is
#    free.
# 4) Two brand new "ruintown" (灰烬镇, lit. "ashtown") rows are appended at
#    the bottom as rows 26-27, continuing the Id sequence (42010016,
#    42010017).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: remove the stray duplicate "portal" row (old row 24) ---
$ws.Rows.Item(24).Delete()

# --- Step 2: insert a fresh row at 11 or the new portal entry, pushing ---
# --- the former rows 11-23 down to 12-24                              ---
$ws.Rows.Item(11).Insert()
$ws.Rows.Item(11).ClearFormats()

$ws.Range("A11").Value = 42010001
$ws.Range("B11").Value = "传送门"
$ws.Range("C11").Value = 0
$ws.Range("E11").Value = "portal"
$ws.Range("F11").Value = "portal"
$ws.Range("G11").Value = "portal"

# --- Step 3: renumber the old final row (now sitting at row 25 again, ---
# --- since the delete above and the insert above cancel out) ---
$ws.Range("A25").Value = 42010015

# --- Step 4: append the two new ruintown rows ---
$ws.Range("A26").Value = 42010016
$ws.Range("B26").Value = "灰烬镇I"
$ws.Range("C26").Value = 0
$ws.Range("E26").Value = "ruintown1"
$ws.Range("F26").Value = "ruintown"
$ws.Range("G26").Value = "ruintown1"
$ws.Range("L26").Value = 50

$ws.Range("A27").Value = 42010017
$ws.Range("B27").Value = "灰烬镇II"
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = "ruintown2"
$ws.Range("F27").Value = "ruintown"
$ws.Range("G27").Value = "ruintown2"
$ws.Range("I27").Value = 43000031
$ws.Range("J27").Value = "oneline"
$ws.Range("K27").Value = 200
$ws.Range("L27").Value = 100
$ws.Range("O27").Value = 100
$ws.Range("Q27").Value = 23000102
$ws.Range("T27").Value = 100

# --- Step 5: grow the table (ListObject) so it covers the new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A3:Y27"))

# --- Step 6: move the selection the way the author left it ---
$ws.Range("A23").Select()
